$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Maracuyá - date moved to 44417, values swapped with row 3)
$ws.Range("D2").Value = 44417
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 26000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 26000
$ws.Range("R2").Value = "Perú"
$ws.Range("S2").Value = 1444

# Row 3 updates (date moved to 44235, values swapped with row 2)
$ws.Range("D3").Value = 44235
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 42000
$ws.Range("O3").Value = 42000
$ws.Range("P3").Value = 42000
$ws.Range("R3").Value = "Región de Arica y Parinacota"
$ws.Range("S3").Value = 2333
